$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2025-10-21 Tuesday" "2025-10-22 Wednesday"

Replace-Text "369÷7=" "649÷5="
Replace-Text "633÷5=" "278÷8="
Replace-Text "406÷6=" "565÷8="
Replace-Text "423÷6=" "755÷2="
Replace-Text "177÷8=" "663÷5="
Replace-Text "674÷2=" "443÷6="
Replace-Text "553÷2=" "897÷4="
Replace-Text "373÷7=" "939÷7="
Replace-Text "312÷6=" "296÷4="
Replace-Text "903÷3=" "257÷9="
Replace-Text "404÷7=" "696÷2="
Replace-Text "641÷9=" "491÷7="
Replace-Text "433÷5=" "994÷6="
Replace-Text "786÷9=" "386÷3="
Replace-Text "588÷7=" "841÷3="
Replace-Text "519÷9=" "408÷8="
Replace-Text "157÷2=" "228÷2="
Replace-Text "320÷6=" "679÷6="
Replace-Text "706÷8=" "252÷6="
Replace-Text "716÷6=" "622÷7="
Replace-Text "136÷6=" "432÷9="
Replace-Text "577÷4=" "541÷5="
Replace-Text "455÷4=" "836÷7="
Replace-Text "302÷2=" "469÷9="
Replace-Text "948÷8=" "571÷8="

Write-Output "Done"
